# Daily "Updated symbol list" refresh (GitHub Actions bot).
# Price (column D) and Volume(1h) (column E) cells are stored as plain
# text in this sheet, so each cell is formatted as Text ("@") before its
# new value is written - this keeps the refreshed figures as literal text
# (e.g. "283.22", "2.04%") instead of Excel auto-converting them to
# numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
}

Set-TextValue "D2" "283.22"
Set-TextValue "E2" "2.04%"
Set-TextValue "D3" "28.55"
Set-TextValue "E3" "4.55%"
Set-TextValue "D4" "5.071"
Set-TextValue "E4" "3.79%"
Set-TextValue "D5" "0.06485"
Set-TextValue "E5" "1.02%"
Set-TextValue "D6" "7.221"
Set-TextValue "E6" "3.60%"
Set-TextValue "D7" "1.422"
Set-TextValue "E7" "20.29%"
Set-TextValue "D8" "0.9104"
Set-TextValue "E8" "3.07%"
Set-TextValue "D9" "0.1546"
Set-TextValue "E9" "0.33%"
Set-TextValue "D10" "0.06496"
Set-TextValue "E10" "26.66%"
Set-TextValue "D11" "0.07571"
Set-TextValue "E11" "2.15%"
Set-TextValue "D12" "0.02755"
Set-TextValue "E12" "-4.45%"
Set-TextValue "D13" "0.08948"
Set-TextValue "E13" "-0.27%"
Set-TextValue "D14" "0.001597"
Set-TextValue "E14" "1.92%"
Set-TextValue "D15" "0.0006365"
Set-TextValue "E15" "-0.36%"
Set-TextValue "D16" "0.006104"
Set-TextValue "E16" "-0.98%"
Set-TextValue "D17" "3.451"
Set-TextValue "E17" "-0.78%"
Set-TextValue "D18" "3.369"
Set-TextValue "E18" "1.60%"
Set-TextValue "E19" "-1.41%"
Set-TextValue "E20" "1.19%"
Set-TextValue "D21" "0.1342"
Set-TextValue "E21" "0.81%"
Set-TextValue "D22" "3.983"
Set-TextValue "E22" "2.10%"
Set-TextValue "E23" "4.22%"
Set-TextValue "D24" "0.04453"
Set-TextValue "E24" "0.89%"
Set-TextValue "D25" "0.001184"
Set-TextValue "E25" "0.62%"
Set-TextValue "D26" "0.004458"
Set-TextValue "E26" "15.28%"
Set-TextValue "D28" "0.0001200"
Set-TextValue "E28" "1.66%"
Set-TextValue "E29" "-15.63%"
Set-TextValue "D40" "0.04108"
Set-TextValue "E40" "-0.83%"
Set-TextValue "D41" "0.006639"
Set-TextValue "E41" "-2.15%"
Set-TextValue "D42" "0.1229"
Set-TextValue "E42" "4.70%"
Set-TextValue "D43" "0.002050"
Set-TextValue "D44" "0.01239"
Set-TextValue "E44" "8.29%"
Set-TextValue "D45" "0.00005403"
Set-TextValue "E45" "1.63%"
Set-TextValue "E46" "-0.11%"
Set-TextValue "D47" "2.041"
Set-TextValue "E47" "20.97%"
